$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all crops")
$ws.Name = "All crops"
